# Auto-generated cell updates applying the Halicarnassus_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 345.2857
$ws.Range("I2").Value = 124.666664
$ws.Range("K2").Value = 124.666664
$ws.Range("M2").Value = -11.666664
$ws.Range("H40").Value = 4717.615
$ws.Range("I40").Value = 3824.4443
$ws.Range("J40").Value = 6727.25
$ws.Range("K40").Value = 3824.4443
$ws.Range("L40").Value = 6727.25
$ws.Range("M40").Value = -3649.4443
$ws.Range("N40").Value = -7077.25
$ws.Range("H48").Value = 5666.5
$ws.Range("J48").Value = 5666.5
$ws.Range("L48").Value = 16999.5
$ws.Range("N48").Value = -17583.5
$ws.Range("H51").Value = 3500
$ws.Range("J51").Value = 3500
$ws.Range("L51").Value = 3500
$ws.Range("N51").Value = -4468
$ws.Range("H56").Value = 5666.5
$ws.Range("J56").Value = 5666.5
$ws.Range("L56").Value = 16999.5
$ws.Range("N56").Value = -18067.5
$ws.Range("H62").Value = 10598.6
$ws.Range("I62").Value = 6994.5
$ws.Range("J62").Value = 11499.625
$ws.Range("K62").Value = 6994.5
$ws.Range("L62").Value = 11499.625
$ws.Range("M62").Value = -6370.5
$ws.Range("N62").Value = -12747.625
$ws.Range("H65").Value = 10598.6
$ws.Range("I65").Value = 6994.5
$ws.Range("J65").Value = 11499.625
$ws.Range("K65").Value = 34972.5
$ws.Range("L65").Value = 57498.125
$ws.Range("M65").Value = -31852.5
$ws.Range("N65").Value = -63738.125
$ws.Range("H76").Value = 3700.2856
$ws.Range("I76").Value = 3180.4
$ws.Range("K76").Value = 3180.4
$ws.Range("M76").Value = -2865.4
$ws.Range("H79").Value = 3700.2856
$ws.Range("I79").Value = 3180.4
$ws.Range("K79").Value = 3180.4
$ws.Range("M79").Value = -2088.4
$ws.Range("H97").Value = 861.75
$ws.Range("J97").Value = 861.75
$ws.Range("L97").Value = 2585.25
$ws.Range("N97").Value = -3577.25
$ws.Range("H106").Value = 7124.5
$ws.Range("I106").Value = 7124.5
$ws.Range("K106").Value = 7124.5
$ws.Range("M106").Value = -6493.5
$ws.Range("H107").Value = 186.66667
$ws.Range("I107").Value = 124.2
$ws.Range("K107").Value = 124.2
$ws.Range("M107").Value = 1795.8
$ws.Range("H141").Value = 1915.8889
$ws.Range("I141").Value = 1517.875
$ws.Range("K141").Value = 4553.625
$ws.Range("M141").Value = 626.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8471.895
$ws.Range("I32").Value = 6586.2354
$ws.Range("K32").Value = 6586.2354
$ws.Range("M32").Value = -6299.2354
$ws.Range("H62").Value = 35833
$ws.Range("J62").Value = 35833
$ws.Range("L62").Value = 35833
$ws.Range("N62").Value = -37081
$ws.Range("H65").Value = 35833
$ws.Range("J65").Value = 35833
$ws.Range("L65").Value = 107499
$ws.Range("N65").Value = -113739
$ws.Range("H76").Value = 18750
$ws.Range("J76").Value = 18750
$ws.Range("L76").Value = 18750
$ws.Range("N76").Value = -19426
$ws.Range("H79").Value = 18750
$ws.Range("J79").Value = 18750
$ws.Range("L79").Value = 18750
$ws.Range("N79").Value = -21090
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H94").Value = 30500
$ws.Range("J94").Value = 30500
$ws.Range("L94").Value = 30500
$ws.Range("N94").Value = -32302
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
$ws.Range("H96").Value = 5031473.5
$ws.Range("J96").Value = 5031473.5
$ws.Range("L96").Value = 5031473.5
$ws.Range("N96").Value = -5036965.5
$ws.Range("H130").Value = 69397.8
$ws.Range("J130").Value = 69397.8
$ws.Range("L130").Value = 69397.8
$ws.Range("N130").Value = -79437.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7486.5713
$ws.Range("I20").Value = 9474.5
$ws.Range("J20").Value = 4836
$ws.Range("K20").Value = 9474.5
$ws.Range("L20").Value = 4836
$ws.Range("M20").Value = -9227.5
$ws.Range("N20").Value = -5330
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H86").Value = 7900.5454
$ws.Range("I86").Value = 2995
$ws.Range("J86").Value = 8990.666999999999
$ws.Range("K86").Value = 2995
$ws.Range("L86").Value = 8990.666999999999
$ws.Range("M86").Value = -1872
$ws.Range("N86").Value = -11236.667
$ws.Range("H89").Value = 7900.5454
$ws.Range("I89").Value = 2995
$ws.Range("J89").Value = 8990.666999999999
$ws.Range("K89").Value = 14975
$ws.Range("L89").Value = 44953.335
$ws.Range("M89").Value = -9359
$ws.Range("N89").Value = -56185.335
$ws.Range("H105").Value = 1805.5625
$ws.Range("I105").Value = 1627.909
$ws.Range("K105").Value = 1627.909
$ws.Range("M105").Value = 119.0909999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1199.75
$ws.Range("J16").Value = 1200
$ws.Range("L16").Value = 1200
$ws.Range("N16").Value = -1774
$ws.Range("H31").Value = 6135.1333
$ws.Range("I31").Value = 2848
$ws.Range("K31").Value = 2848
$ws.Range("M31").Value = -2553
$ws.Range("H34").Value = 6135.1333
$ws.Range("I34").Value = 2848
$ws.Range("K34").Value = 2848
$ws.Range("M34").Value = -2646
$ws.Range("H35").Value = 186.25
$ws.Range("I35").Value = 186.25
$ws.Range("K35").Value = 186.25
$ws.Range("M35").Value = 107.75
$ws.Range("H56").Value = 15531
$ws.Range("I56").Value = 15531
$ws.Range("K56").Value = 15531
$ws.Range("M56").Value = -14686
$ws.Range("H88").Value = 10208
$ws.Range("J88").Value = 11124.75
$ws.Range("L88").Value = 11124.75
$ws.Range("N88").Value = -11936.75
$ws.Range("H91").Value = 10208
$ws.Range("J91").Value = 11124.75
$ws.Range("L91").Value = 11124.75
$ws.Range("N91").Value = -13932.75
$ws.Range("H113").Value = 1199.75
$ws.Range("J113").Value = 1200
$ws.Range("L113").Value = 1200
$ws.Range("N113").Value = -5540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1750
$ws.Range("I20").Value = 2000
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 4500
$ws.Range("M20").Value = -5773
$ws.Range("N20").Value = -4954
$ws.Range("H21").Value = 450
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null
$ws.Range("H88").Value = 14999
$ws.Range("J88").Value = 14999
$ws.Range("L88").Value = 44997
$ws.Range("N88").Value = -45853
$ws.Range("H91").Value = 14999
$ws.Range("J91").Value = 14999
$ws.Range("L91").Value = 44997
$ws.Range("N91").Value = -47961

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6278.7
$ws.Range("I70").Value = 5348.75
$ws.Range("J70").Value = 9998.5
$ws.Range("K70").Value = 5348.75
$ws.Range("L70").Value = 9998.5
$ws.Range("M70").Value = -5078.75
$ws.Range("N70").Value = -10538.5
$ws.Range("H73").Value = 6278.7
$ws.Range("I73").Value = 5348.75
$ws.Range("J73").Value = 9998.5
$ws.Range("K73").Value = 5348.75
$ws.Range("L73").Value = 9998.5
$ws.Range("M73").Value = -4412.75
$ws.Range("N73").Value = -11870.5
$ws.Range("H102").Value = 3566.6365
$ws.Range("I102").Value = 3523.3
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 3523.3
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -1901.3
$ws.Range("N102").Value = -7244
$ws.Range("H107").Value = 450
$ws.Range("I107").Value = 300
$ws.Range("K107").Value = 300
$ws.Range("M107").Value = 1620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9621.75
$ws.Range("I7").Value = 8499
$ws.Range("K7").Value = 8499
$ws.Range("M7").Value = -8387
$ws.Range("H22").Value = 1655
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705
$ws.Range("H27").Value = 1655
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893
$ws.Range("H126").Value = 9621.75
$ws.Range("I126").Value = 8499
$ws.Range("K126").Value = 25497
$ws.Range("M126").Value = -23027

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 49989
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 49989
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H100").Value = 1453.8
$ws.Range("I100").Value = 1667.25
$ws.Range("K100").Value = 3334.5
$ws.Range("M100").Value = -2793.5
$ws.Range("H107").Value = 521
$ws.Range("I107").Value = 463.08334
$ws.Range("K107").Value = 1389.25002
$ws.Range("M107").Value = 530.7499800000001
$ws.Range("H132").Value = 2536.4
$ws.Range("I132").Value = 1945.8462
$ws.Range("K132").Value = 5837.5386
$ws.Range("M132").Value = -3307.5386

Write-Host "Applied $(53) cell updates across 8 sheets"